$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 88.3125  # ALC!H33: 82.111115 -> 88.3125
$ws.Cells.Item(33, 9).Value = 88.3125  # ALC!I33: 82.111115 -> 88.3125
$ws.Cells.Item(33, 11).Value = 88.3125  # ALC!K33: 82.111115 -> 88.3125
$ws.Cells.Item(33, 13).Value = 140.6875  # ALC!M33: 146.888885 -> 140.6875

$ws.Cells.Item(112, 8).Value = 1648.579  # ALC!H112: 1566.5 -> 1648.579
$ws.Cells.Item(112, 10).Value = 1648.579  # ALC!J112: 1566.5 -> 1648.579
$ws.Cells.Item(112, 12).Value = 4945.737  # ALC!L112: 4699.5 -> 4945.737
$ws.Cells.Item(112, 14).Value = -7161.737  # ALC!N112: -6915.5 -> -7161.737

$ws.Cells.Item(113, 8).Value = 33336960  # ALC!H113: 37041092 -> 33336960
$ws.Cells.Item(113, 9).Value = 47621776  # ALC!I113: 66669330 -> 47621776
$ws.Cells.Item(113, 10).Value = 5728.6665  # ALC!J113: 5798 -> 5728.6665
$ws.Cells.Item(113, 11).Value = 47621776  # ALC!K113: 66669330 -> 47621776
$ws.Cells.Item(113, 12).Value = 5728.6665  # ALC!L113: 5798 -> 5728.6665
$ws.Cells.Item(113, 13).Value = -47618522  # ALC!M113: -66666076 -> -47618522
$ws.Cells.Item(113, 14).Value = -12236.6665  # ALC!N113: -12306 -> -12236.6665

$ws.Cells.Item(121, 8).Value = 4997.5  # ALC!H121: 5000 -> 4997.5
$ws.Cells.Item(121, 10).Value = 4997.5  # ALC!J121: 5000 -> 4997.5
$ws.Cells.Item(121, 12).Value = 14992.5  # ALC!L121: 15000 -> 14992.5
$ws.Cells.Item(121, 14).Value = -18486.5  # ALC!N121: -18494 -> -18486.5

$ws.Cells.Item(129, 8).Value = 1362.9546  # ALC!H129: 1316.6957 -> 1362.9546
$ws.Cells.Item(129, 9).Value = 491.85715  # ALC!I129: 479 -> 491.85715
$ws.Cells.Item(129, 11).Value = 1475.57145  # ALC!K129: 1437 -> 1475.57145
$ws.Cells.Item(129, 13).Value = 3524.42855  # ALC!M129: 3563 -> 3524.42855

$ws.Cells.Item(135, 8).Value = 2407.4167  # ALC!H135: 2308.9 -> 2407.4167
$ws.Cells.Item(135, 9).Value = 1543.5555  # ALC!I135: 1636.5 -> 1543.5555
$ws.Cells.Item(135, 10).Value = 4999  # ALC!J135: 4998.5 -> 4999
$ws.Cells.Item(135, 11).Value = 13891.9995  # ALC!K135: 14728.5 -> 13891.9995
$ws.Cells.Item(135, 12).Value = 44991  # ALC!L135: 44986.5 -> 44991
$ws.Cells.Item(135, 13).Value = -11356.9995  # ALC!M135: -12193.5 -> -11356.9995
$ws.Cells.Item(135, 14).Value = -50061  # ALC!N135: -50056.5 -> -50061

$ws.Cells.Item(141, 8).Value = 2536.9375  # ALC!H141: 2319.5334 -> 2536.9375
$ws.Cells.Item(141, 9).Value = 2306.0667  # ALC!I141: 2319.5334 -> 2306.0667
$ws.Cells.Item(141, 10).Value = 6000  # ALC!J141: 0 -> 6000
$ws.Cells.Item(141, 11).Value = 6918.2001  # ALC!K141: 6958.600199999999 -> 6918.2001
$ws.Cells.Item(141, 12).Value = 18000  # ALC!L141: 0 -> 18000
$ws.Cells.Item(141, 13).Value = -1738.2001  # ALC!M141: -1778.600199999999 -> -1738.2001
$ws.Cells.Item(141, 14).Value = -28360  # ALC!N141: <MISSING> -> -28360

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 3055.111  # ARM!H61: 3102.4707 -> 3055.111
$ws.Cells.Item(61, 9).Value = 2453.9092  # ARM!I61: 2474.3 -> 2453.9092
$ws.Cells.Item(61, 11).Value = 2453.9092  # ARM!K61: 2474.3 -> 2453.9092
$ws.Cells.Item(61, 13).Value = -2241.9092  # ARM!M61: -2262.3 -> -2241.9092

$ws.Cells.Item(74, 8).Value = 1598.7333  # ARM!H74: 1601.7667 -> 1598.7333
$ws.Cells.Item(74, 9).Value = 1072.55  # ARM!I74: 1059.0952 -> 1072.55
$ws.Cells.Item(74, 10).Value = 2651.1  # ARM!J74: 2868 -> 2651.1
$ws.Cells.Item(74, 11).Value = 1072.55  # ARM!K74: 1059.0952 -> 1072.55
$ws.Cells.Item(74, 12).Value = 2651.1  # ARM!L74: 2868 -> 2651.1
$ws.Cells.Item(74, 13).Value = -198.55  # ARM!M74: -185.0952 -> -198.55
$ws.Cells.Item(74, 14).Value = -4399.1  # ARM!N74: -4616 -> -4399.1

$ws.Cells.Item(77, 8).Value = 1598.7333  # ARM!H77: 1601.7667 -> 1598.7333
$ws.Cells.Item(77, 9).Value = 1072.55  # ARM!I77: 1059.0952 -> 1072.55
$ws.Cells.Item(77, 10).Value = 2651.1  # ARM!J77: 2868 -> 2651.1
$ws.Cells.Item(77, 11).Value = 5362.75  # ARM!K77: 5295.476 -> 5362.75
$ws.Cells.Item(77, 12).Value = 13255.5  # ARM!L77: 14340 -> 13255.5
$ws.Cells.Item(77, 13).Value = -994.75  # ARM!M77: -927.4759999999997 -> -994.75
$ws.Cells.Item(77, 14).Value = -21991.5  # ARM!N77: -23076 -> -21991.5

$ws.Cells.Item(97, 8).Value = 671.875  # ARM!H97: 650.2353000000001 -> 671.875
$ws.Cells.Item(97, 9).Value = 567  # ARM!I97: 546.7692 -> 567
$ws.Cells.Item(97, 11).Value = 567  # ARM!K97: 546.7692 -> 567
$ws.Cells.Item(97, 13).Value = -71  # ARM!M97: -50.76919999999996 -> -71

$ws.Cells.Item(102, 8).Value = 2652.2307  # ARM!H102: 2227.2 -> 2652.2307
$ws.Cells.Item(102, 9).Value = 1886.6666  # ARM!I102: 1690.3125 -> 1886.6666
$ws.Cells.Item(102, 11).Value = 1886.6666  # ARM!K102: 1690.3125 -> 1886.6666
$ws.Cells.Item(102, 13).Value = -264.6666  # ARM!M102: -68.3125 -> -264.6666

$ws.Cells.Item(122, 8).Value = 3375.9285  # ARM!H122: 3575.1538 -> 3375.9285
$ws.Cells.Item(122, 9).Value = 3853.182  # ARM!I122: 4099.3 -> 3853.182
$ws.Cells.Item(122, 10).Value = 1626  # ARM!J122: 1828 -> 1626
$ws.Cells.Item(122, 11).Value = 11559.546  # ARM!K122: 12297.9 -> 11559.546
$ws.Cells.Item(122, 12).Value = 4878  # ARM!L122: 5484 -> 4878
$ws.Cells.Item(122, 13).Value = -9109.545999999998  # ARM!M122: -9847.900000000001 -> -9109.545999999998
$ws.Cells.Item(122, 14).Value = -9778  # ARM!N122: -10384 -> -9778

$ws.Cells.Item(135, 8).Value = 94320  # ARM!H135: 94426.664 -> 94320
$ws.Cells.Item(135, 10).Value = 94320  # ARM!J135: 94426.664 -> 94320
$ws.Cells.Item(135, 12).Value = 94320  # ARM!L135: 94426.664 -> 94320
$ws.Cells.Item(135, 14).Value = -104460  # ARM!N135: -104566.664 -> -104460

$ws.Cells.Item(136, 8).Value = 3055.111  # ARM!H136: 3102.4707 -> 3055.111
$ws.Cells.Item(136, 9).Value = 2453.9092  # ARM!I136: 2474.3 -> 2453.9092
$ws.Cells.Item(136, 11).Value = 7361.7276  # ARM!K136: 7422.900000000001 -> 7361.7276
$ws.Cells.Item(136, 13).Value = -4811.7276  # ARM!M136: -4872.900000000001 -> -4811.7276

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 20837236  # BSM!H20: 22731508 -> 20837236
$ws.Cells.Item(20, 9).Value = 50007020  # BSM!I20: 62508676 -> 50007020
$ws.Cells.Item(20, 10).Value = 1676.0714  # BSM!J20: 1698.0714 -> 1676.0714
$ws.Cells.Item(20, 11).Value = 50007020  # BSM!K20: 62508676 -> 50007020
$ws.Cells.Item(20, 12).Value = 1676.0714  # BSM!L20: 1698.0714 -> 1676.0714
$ws.Cells.Item(20, 13).Value = -50006773  # BSM!M20: -62508429 -> -50006773
$ws.Cells.Item(20, 14).Value = -2170.0714  # BSM!N20: -2192.0714 -> -2170.0714

$ws.Cells.Item(86, 8).Value = 3326.162  # BSM!H86: 3561.639 -> 3326.162
$ws.Cells.Item(86, 9).Value = 3295.3333  # BSM!I86: 3565.08 -> 3295.3333
$ws.Cells.Item(86, 10).Value = 3409.4  # BSM!J86: 3553.818 -> 3409.4
$ws.Cells.Item(86, 11).Value = 3295.3333  # BSM!K86: 3565.08 -> 3295.3333
$ws.Cells.Item(86, 12).Value = 3409.4  # BSM!L86: 3553.818 -> 3409.4
$ws.Cells.Item(86, 13).Value = -2172.3333  # BSM!M86: -2442.08 -> -2172.3333
$ws.Cells.Item(86, 14).Value = -5655.4  # BSM!N86: -5799.818 -> -5655.4

$ws.Cells.Item(89, 8).Value = 3326.162  # BSM!H89: 3561.639 -> 3326.162
$ws.Cells.Item(89, 9).Value = 3295.3333  # BSM!I89: 3565.08 -> 3295.3333
$ws.Cells.Item(89, 10).Value = 3409.4  # BSM!J89: 3553.818 -> 3409.4
$ws.Cells.Item(89, 11).Value = 16476.6665  # BSM!K89: 17825.4 -> 16476.6665
$ws.Cells.Item(89, 12).Value = 17047  # BSM!L89: 17769.09 -> 17047
$ws.Cells.Item(89, 13).Value = -10860.6665  # BSM!M89: -12209.4 -> -10860.6665
$ws.Cells.Item(89, 14).Value = -28279  # BSM!N89: -29001.09 -> -28279

$ws.Cells.Item(99, 8).Value = 2110.5557  # BSM!H99: 2342.0667 -> 2110.5557
$ws.Cells.Item(99, 9).Value = 1489.7273  # BSM!I99: 1512.5454 -> 1489.7273
$ws.Cells.Item(99, 10).Value = 3086.1428  # BSM!J99: 4623.25 -> 3086.1428
$ws.Cells.Item(99, 11).Value = 1489.7273  # BSM!K99: 1512.5454 -> 1489.7273
$ws.Cells.Item(99, 12).Value = 3086.1428  # BSM!L99: 4623.25 -> 3086.1428
$ws.Cells.Item(99, 13).Value = 8.272699999999986  # BSM!M99: -14.54539999999997 -> 8.272699999999986
$ws.Cells.Item(99, 14).Value = -6082.1428  # BSM!N99: -7619.25 -> -6082.1428

$ws.Cells.Item(105, 8).Value = 10835050  # BSM!H105: 11820009 -> 10835050
$ws.Cells.Item(105, 9).Value = 910799  # BSM!I105: 1113088.5 -> 910799
$ws.Cells.Item(105, 11).Value = 910799  # BSM!K105: 1113088.5 -> 910799
$ws.Cells.Item(105, 13).Value = -909052  # BSM!M105: -1111341.5 -> -909052

$ws.Cells.Item(107, 8).Value = 3345738.5  # BSM!H107: 3206345.2 -> 3345738.5
$ws.Cells.Item(107, 9).Value = 5129342  # BSM!I107: 4808777 -> 5129342
$ws.Cells.Item(107, 11).Value = 5129342  # BSM!K107: 4808777 -> 5129342
$ws.Cells.Item(107, 13).Value = -5127422  # BSM!M107: -4806857 -> -5127422

$ws.Cells.Item(134, 8).Value = 3037.9375  # BSM!H134: 3192.3333 -> 3037.9375
$ws.Cells.Item(134, 9).Value = 2244.1428  # BSM!I134: 2482 -> 2244.1428
$ws.Cells.Item(134, 10).Value = 3655.3333  # BSM!J134: 3699.7144 -> 3655.3333
$ws.Cells.Item(134, 11).Value = 6732.428400000001  # BSM!K134: 7446 -> 6732.428400000001
$ws.Cells.Item(134, 12).Value = 10965.9999  # BSM!L134: 11099.1432 -> 10965.9999
$ws.Cells.Item(134, 13).Value = -4197.428400000001  # BSM!M134: -4911 -> -4197.428400000001
$ws.Cells.Item(134, 14).Value = -16035.9999  # BSM!N134: -16169.1432 -> -16035.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(94, 8).Value = 527.8  # CRP!H94: 541.4 -> 527.8
$ws.Cells.Item(94, 9).Value = 453.6  # CRP!I94: 480.8 -> 453.6
$ws.Cells.Item(94, 11).Value = 453.6  # CRP!K94: 480.8 -> 453.6
$ws.Cells.Item(94, 13).Value = -2.600000000000023  # CRP!M94: -29.80000000000001 -> -2.600000000000023

$ws.Cells.Item(99, 8).Value = 4392.1816  # CRP!H99: 4351.9165 -> 4392.1816
$ws.Cells.Item(99, 10).Value = 4946.3335  # CRP!J99: 4798.143 -> 4946.3335
$ws.Cells.Item(99, 12).Value = 4946.3335  # CRP!L99: 4798.143 -> 4946.3335
$ws.Cells.Item(99, 14).Value = -7942.3335  # CRP!N99: -7794.143 -> -7942.3335

$ws.Cells.Item(122, 8).Value = 3066.0557  # CRP!H122: 3011.7058 -> 3066.0557
$ws.Cells.Item(122, 10).Value = 4059.2  # CRP!J122: 4076.5 -> 4059.2
$ws.Cells.Item(122, 12).Value = 12177.6  # CRP!L122: 12229.5 -> 12177.6
$ws.Cells.Item(122, 14).Value = -17077.6  # CRP!N122: -17129.5 -> -17077.6

$ws.Cells.Item(126, 8).Value = 4392.1816  # CRP!H126: 4351.9165 -> 4392.1816
$ws.Cells.Item(126, 10).Value = 4946.3335  # CRP!J126: 4798.143 -> 4946.3335
$ws.Cells.Item(126, 12).Value = 14839.0005  # CRP!L126: 14394.429 -> 14839.0005
$ws.Cells.Item(126, 14).Value = -19779.0005  # CRP!N126: -19334.429 -> -19779.0005

$ws.Cells.Item(134, 8).Value = 4431.643  # CRP!H134: 4376.6 -> 4431.643
$ws.Cells.Item(134, 9).Value = 4808.304  # CRP!I134: 4679.2915 -> 4808.304
$ws.Cells.Item(134, 10).Value = 2699  # CRP!J134: 3165.8333 -> 2699
$ws.Cells.Item(134, 11).Value = 14424.912  # CRP!K134: 14037.8745 -> 14424.912
$ws.Cells.Item(134, 12).Value = 8097  # CRP!L134: 9497.499899999999 -> 8097
$ws.Cells.Item(134, 13).Value = -11889.912  # CRP!M134: -11502.8745 -> -11889.912
$ws.Cells.Item(134, 14).Value = -13167  # CRP!N134: -14567.4999 -> -13167

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 2891.2666  # CUL!H3: 2729.3125 -> 2891.2666
$ws.Cells.Item(3, 9).Value = 2312.0715  # CUL!I3: 2177.9333 -> 2312.0715
$ws.Cells.Item(3, 11).Value = 6936.2145  # CUL!K3: 6533.7999 -> 6936.2145
$ws.Cells.Item(3, 13).Value = -6824.2145  # CUL!M3: -6421.7999 -> -6824.2145

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 752.7692  # GSM!H97: 776.08 -> 752.7692
$ws.Cells.Item(97, 9).Value = 798.1739  # GSM!I97: 826.7273 -> 798.1739
$ws.Cells.Item(97, 11).Value = 798.1739  # GSM!K97: 826.7273 -> 798.1739
$ws.Cells.Item(97, 13).Value = -302.1739  # GSM!M97: -330.7273 -> -302.1739

$ws.Cells.Item(102, 8).Value = 4134.237  # GSM!H102: 4565.8184 -> 4134.237
$ws.Cells.Item(102, 9).Value = 1600.909  # GSM!I102: 1863.5 -> 1600.909
$ws.Cells.Item(102, 11).Value = 1600.909  # GSM!K102: 1863.5 -> 1600.909
$ws.Cells.Item(102, 13).Value = 21.09099999999989  # GSM!M102: -241.5 -> 21.09099999999989

$ws.Cells.Item(105, 8).Value = 77500  # GSM!H105: 0 -> 77500
$ws.Cells.Item(105, 10).Value = 77500  # GSM!J105: 0 -> 77500
$ws.Cells.Item(105, 12).Value = 77500  # GSM!L105: 0 -> 77500
$ws.Cells.Item(105, 14).Value = -84488  # GSM!N105: <MISSING> -> -84488

$ws.Cells.Item(132, 8).Value = 2492.5833  # GSM!H132: 2394.3572 -> 2492.5833
$ws.Cells.Item(132, 9).Value = 2337.6667  # GSM!I132: 2124.6 -> 2337.6667
$ws.Cells.Item(132, 11).Value = 7013.000100000001  # GSM!K132: 6373.799999999999 -> 7013.000100000001
$ws.Cells.Item(132, 13).Value = -4483.000100000001  # GSM!M132: -3843.799999999999 -> -4483.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3086.9092  # LTW!H7: 3078.5833 -> 3086.9092
$ws.Cells.Item(7, 9).Value = 2990.4  # LTW!I7: 2989.8333 -> 2990.4
$ws.Cells.Item(7, 11).Value = 2990.4  # LTW!K7: 2989.8333 -> 2990.4
$ws.Cells.Item(7, 13).Value = -2878.4  # LTW!M7: -2877.8333 -> -2878.4

$ws.Cells.Item(40, 8).Value = 33636.547  # LTW!H40: 46285.57 -> 33636.547
$ws.Cells.Item(40, 9).Value = 33636.547  # LTW!I40: 46285.57 -> 33636.547
$ws.Cells.Item(40, 11).Value = 33636.547  # LTW!K40: 46285.57 -> 33636.547
$ws.Cells.Item(40, 13).Value = -33500.547  # LTW!M40: -46149.57 -> -33500.547

$ws.Cells.Item(93, 8).Value = 2134.9412  # LTW!H93: 2213 -> 2134.9412
$ws.Cells.Item(93, 9).Value = 2368.6316  # LTW!I93: 2508.3157 -> 2368.6316
$ws.Cells.Item(93, 11).Value = 2368.6316  # LTW!K93: 2508.3157 -> 2368.6316
$ws.Cells.Item(93, 13).Value = -1120.6316  # LTW!M93: -1260.3157 -> -1120.6316

$ws.Cells.Item(122, 8).Value = 4688.5  # LTW!H122: 4728.9 -> 4688.5
$ws.Cells.Item(122, 9).Value = 2333  # LTW!I122: 2467.6667 -> 2333
$ws.Cells.Item(122, 11).Value = 6999  # LTW!K122: 7403.000100000001 -> 6999
$ws.Cells.Item(122, 13).Value = -4549  # LTW!M122: -4953.000100000001 -> -4549

$ws.Cells.Item(126, 8).Value = 3086.9092  # LTW!H126: 3078.5833 -> 3086.9092
$ws.Cells.Item(126, 9).Value = 2990.4  # LTW!I126: 2989.8333 -> 2990.4
$ws.Cells.Item(126, 11).Value = 8971.200000000001  # LTW!K126: 8969.499899999999 -> 8971.200000000001
$ws.Cells.Item(126, 13).Value = -6501.200000000001  # LTW!M126: -6499.499899999999 -> -6501.200000000001

$ws.Cells.Item(132, 8).Value = 4584.5  # LTW!H132: 5082.4287 -> 4584.5
$ws.Cells.Item(132, 9).Value = 4726.4  # LTW!I132: 5633.25 -> 4726.4
$ws.Cells.Item(132, 11).Value = 14179.2  # LTW!K132: 16899.75 -> 14179.2
$ws.Cells.Item(132, 13).Value = -11649.2  # LTW!M132: -14369.75 -> -11649.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(39, 8).Value = 48666.332  # WVR!H39: 25624.25 -> 48666.332
$ws.Cells.Item(39, 9).Value = 0  # WVR!I39: 3500 -> 0
$ws.Cells.Item(39, 10).Value = 48666.332  # WVR!J39: 32999 -> 48666.332
$ws.Cells.Item(39, 11).Value = 0  # WVR!K39: 3500 -> 0
$ws.Cells.Item(39, 12).Value = 48666.332  # WVR!L39: 32999 -> 48666.332
$ws.Cells.Item(39, 13).ClearContents()  # WVR!M39: -3087 -> (removed)
$ws.Cells.Item(39, 14).Value = -49492.332  # WVR!N39: -33825 -> -49492.332

$ws.Cells.Item(81, 8).Value = 3632.682  # WVR!H81: 3489.9565 -> 3632.682
$ws.Cells.Item(81, 9).Value = 4080.3635  # WVR!I81: 3769.5 -> 4080.3635
$ws.Cells.Item(81, 11).Value = 8160.727  # WVR!K81: 7539 -> 8160.727
$ws.Cells.Item(81, 13).Value = -7099.727  # WVR!M81: -6478 -> -7099.727

$ws.Cells.Item(84, 8).Value = 3632.682  # WVR!H84: 3489.9565 -> 3632.682
$ws.Cells.Item(84, 9).Value = 4080.3635  # WVR!I84: 3769.5 -> 4080.3635
$ws.Cells.Item(84, 11).Value = 40803.635  # WVR!K84: 37695 -> 40803.635
$ws.Cells.Item(84, 13).Value = -35499.635  # WVR!M84: -32391 -> -35499.635

$ws.Cells.Item(100, 8).Value = 142858270  # WVR!H100: 100000936 -> 142858270
$ws.Cells.Item(100, 9).Value = 1453.25  # WVR!I100: 1041 -> 1453.25
$ws.Cells.Item(100, 11).Value = 2906.5  # WVR!K100: 2082 -> 2906.5
$ws.Cells.Item(100, 13).Value = -2365.5  # WVR!M100: -1541 -> -2365.5

$ws.Cells.Item(107, 8).Value = 615.5  # WVR!H107: 628.41174 -> 615.5
$ws.Cells.Item(107, 10).Value = 642.5  # WVR!J107: 724.6667 -> 642.5
$ws.Cells.Item(107, 12).Value = 1927.5  # WVR!L107: 2174.0001 -> 1927.5
$ws.Cells.Item(107, 14).Value = -5767.5  # WVR!N107: -6014.0001 -> -5767.5

$ws.Cells.Item(122, 8).Value = 13891010  # WVR!H122: 16668872 -> 13891010
$ws.Cells.Item(122, 9).Value = 2265.3333  # WVR!I122: 2352.3076 -> 2265.3333
$ws.Cells.Item(122, 10).Value = 83334730  # WVR!J122: 125001250 -> 83334730
$ws.Cells.Item(122, 11).Value = 6795.999899999999  # WVR!K122: 7056.9228 -> 6795.999899999999
$ws.Cells.Item(122, 12).Value = 250004190  # WVR!L122: 375003750 -> 250004190
$ws.Cells.Item(122, 13).Value = -4345.999899999999  # WVR!M122: -4606.9228 -> -4345.999899999999
$ws.Cells.Item(122, 14).Value = -250009090  # WVR!N122: -375008650 -> -250009090

$ws.Cells.Item(126, 8).Value = 11679.417  # WVR!H126: 11165.538 -> 11679.417
$ws.Cells.Item(126, 10).Value = 4499.5  # WVR!J126: 4666 -> 4499.5
$ws.Cells.Item(126, 12).Value = 13498.5  # WVR!L126: 13998 -> 13498.5
$ws.Cells.Item(126, 14).Value = -18438.5  # WVR!N126: -18938 -> -18438.5

$ws.Cells.Item(132, 8).Value = 3690.7334  # WVR!H132: 3221.0952 -> 3690.7334
$ws.Cells.Item(132, 9).Value = 4235.6665  # WVR!I132: 3309.875 -> 4235.6665
$ws.Cells.Item(132, 10).Value = 2873.3333  # WVR!J132: 2937 -> 2873.3333
$ws.Cells.Item(132, 11).Value = 12706.9995  # WVR!K132: 9929.625 -> 12706.9995
$ws.Cells.Item(132, 12).Value = 8619.999899999999  # WVR!L132: 8811 -> 8619.999899999999
$ws.Cells.Item(132, 13).Value = -10176.9995  # WVR!M132: -7399.625 -> -10176.9995
$ws.Cells.Item(132, 14).Value = -13679.9999  # WVR!N132: -13871 -> -13679.9999

$ws.Cells.Item(136, 8).Value = 3522.923  # WVR!H136: 2844.9375 -> 3522.923
$ws.Cells.Item(136, 9).Value = 3522.923  # WVR!I136: 3192.4285 -> 3522.923
$ws.Cells.Item(136, 10).Value = 0  # WVR!J136: 412.5 -> 0
$ws.Cells.Item(136, 11).Value = 10568.769  # WVR!K136: 9577.2855 -> 10568.769
$ws.Cells.Item(136, 12).Value = 0  # WVR!L136: 1237.5 -> 0
$ws.Cells.Item(136, 13).Value = -8018.769  # WVR!M136: -7027.2855 -> -8018.769
$ws.Cells.Item(136, 14).ClearContents()  # WVR!N136: -6337.5 -> (removed)
